$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 2935.2856
$ws.Range("I113").Value = 3500
$ws.Range("K113").Value = 3500
$ws.Range("M113").Value = -246
# Row 123
$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800
# Row 132
$ws.Range("H132").Value = 3760894.2
$ws.Range("I132").Value = 4083056.8
$ws.Range("J132").Value = 2333.3333
$ws.Range("K132").Value = 12249170.4
$ws.Range("L132").Value = 6999.999899999999
$ws.Range("M132").Value = -12246640.4
$ws.Range("N132").Value = -12059.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -623
$ws.Range("N45").Value = -1754
# Row 74
$ws.Range("H74").Value = 1712.5
$ws.Range("I74").Value = 1483.3334
$ws.Range("K74").Value = 1483.3334
$ws.Range("M74").Value = -609.3334
# Row 77
$ws.Range("H77").Value = 1712.5
$ws.Range("I77").Value = 1483.3334
$ws.Range("K77").Value = 7416.666999999999
$ws.Range("M77").Value = -3048.666999999999
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0
# Row 132
$ws.Range("H132").Value = 4511.4688
$ws.Range("I132").Value = 4475.115
$ws.Range("K132").Value = 13425.345
$ws.Range("M132").Value = -10895.345

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 20004938
$ws.Range("I31").Value = 5485.778
$ws.Range("J31").Value = 200000000
$ws.Range("K31").Value = 5485.778
$ws.Range("L31").Value = 200000000
$ws.Range("M31").Value = -5190.778
$ws.Range("N31").Value = -200000590
# Row 34
$ws.Range("H34").Value = 20004938
$ws.Range("I34").Value = 5485.778
$ws.Range("J34").Value = 200000000
$ws.Range("K34").Value = 5485.778
$ws.Range("L34").Value = 200000000
$ws.Range("M34").Value = -5283.778
$ws.Range("N34").Value = -200000404
# Row 99
$ws.Range("H99").Value = 1615.7742
$ws.Range("I99").Value = 1167.8422
$ws.Range("J99").Value = 2325
$ws.Range("K99").Value = 1167.8422
$ws.Range("L99").Value = 2325
$ws.Range("M99").Value = 330.1578
$ws.Range("N99").Value = -5321
# Row 107
$ws.Range("H107").Value = 490.44
$ws.Range("J107").Value = 444.125
$ws.Range("L107").Value = 444.125
$ws.Range("N107").Value = -4284.125
# Row 126
$ws.Range("H126").Value = 1615.7742
$ws.Range("I126").Value = 1167.8422
$ws.Range("J126").Value = 2325
$ws.Range("K126").Value = 3503.5266
$ws.Range("L126").Value = 6975
$ws.Range("M126").Value = -1033.5266
$ws.Range("N126").Value = -11915
# Row 129
$ws.Range("H129").Value = 45699.4
$ws.Range("J129").Value = 45699.4
$ws.Range("L129").Value = 45699.4
$ws.Range("N129").Value = -55699.4
# Row 132
$ws.Range("H132").Value = 2189.36
$ws.Range("I132").Value = 1511.238
$ws.Range("K132").Value = 4533.714
$ws.Range("M132").Value = -2003.714
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 8700
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 9933.333000000001
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 9933.333000000001
$ws.Range("M43").Value = -4849
$ws.Range("N43").Value = -10235.333
# Row 80
$ws.Range("H80").Value = 6593.125
$ws.Range("I80").Value = 4862.5
$ws.Range("K80").Value = 4862.5
$ws.Range("M80").Value = -3864.5
# Row 83
$ws.Range("H83").Value = 6593.125
$ws.Range("I83").Value = 4862.5
$ws.Range("K83").Value = 24312.5
$ws.Range("M83").Value = -19320.5
# Row 102
$ws.Range("H102").Value = 1482.3334
$ws.Range("I102").Value = 1374.1765
$ws.Range("J102").Value = 1745
$ws.Range("K102").Value = 1374.1765
$ws.Range("L102").Value = 1745
$ws.Range("M102").Value = 247.8235
$ws.Range("N102").Value = -4989
# Row 107
$ws.Range("H107").Value = 525.7931
$ws.Range("I107").Value = 453.69565
$ws.Range("J107").Value = 802.1667
$ws.Range("K107").Value = 453.69565
$ws.Range("L107").Value = 802.1667
$ws.Range("M107").Value = 1466.30435
$ws.Range("N107").Value = -4642.1667
# Row 113
$ws.Range("H113").Value = 25000992
$ws.Range("J113").Value = 1219
$ws.Range("L113").Value = 1219
$ws.Range("N113").Value = -5559

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2930.3845
$ws.Range("I7").Value = 2313.5715
$ws.Range("J7").Value = 3650
$ws.Range("K7").Value = 2313.5715
$ws.Range("L7").Value = 3650
$ws.Range("M7").Value = -2201.5715
$ws.Range("N7").Value = -3874
# Row 22
$ws.Range("H22").Value = 778.8570999999999
$ws.Range("I22").Value = 187.5
$ws.Range("J22").Value = 1567.3334
$ws.Range("K22").Value = 187.5
$ws.Range("L22").Value = 1567.3334
$ws.Range("M22").Value = 107.5
$ws.Range("N22").Value = -2157.3334
# Row 27
$ws.Range("H27").Value = 778.8570999999999
$ws.Range("I27").Value = 187.5
$ws.Range("J27").Value = 1567.3334
$ws.Range("K27").Value = 187.5
$ws.Range("L27").Value = 1567.3334
$ws.Range("M27").Value = -80.5
$ws.Range("N27").Value = -1781.3334
# Row 40
$ws.Range("H40").Value = 2694.5715
$ws.Range("I40").Value = 2268.6667
$ws.Range("K40").Value = 2268.6667
$ws.Range("M40").Value = -2132.6667
# Row 46
$ws.Range("H46").Value = 1784.8695
$ws.Range("I46").Value = 1166.8334
$ws.Range("J46").Value = 2003
$ws.Range("K46").Value = 1166.8334
$ws.Range("L46").Value = 2003
$ws.Range("M46").Value = -978.8334
$ws.Range("N46").Value = -2379
# Row 64
$ws.Range("H64").Value = 29800
$ws.Range("J64").Value = 29800
$ws.Range("L64").Value = 29800
$ws.Range("N64").Value = -30250
# Row 67
$ws.Range("H67").Value = 29800
$ws.Range("J67").Value = 29800
$ws.Range("L67").Value = 29800
$ws.Range("N67").Value = -31360
# Row 126
$ws.Range("H126").Value = 2930.3845
$ws.Range("I126").Value = 2313.5715
$ws.Range("J126").Value = 3650
$ws.Range("K126").Value = 6940.7145
$ws.Range("L126").Value = 10950
$ws.Range("M126").Value = -4470.7145
$ws.Range("N126").Value = -15890
# Row 128
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 2800
$ws.Range("I10").Value = 2800
$ws.Range("K10").Value = 2800
$ws.Range("M10").Value = -2631
# Row 124
$ws.Range("H124").Value = 39714.5
$ws.Range("J124").Value = 39714.5
$ws.Range("L124").Value = 39714.5
$ws.Range("N124").Value = -49534.5
